$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.073035276108662
$ws.Cells.Item(2, 4).Value = 1.07616131307801
$ws.Cells.Item(2, 5).Value = 1.074059216939659
$ws.Cells.Item(2, 6).Value = 1.086184276494918
$ws.Cells.Item(2, 9).Value = 1.058092633451352
$ws.Cells.Item(2, 10).Value = 1.077951812766228
$ws.Cells.Item(2, 11).Value = 1.078845974784472
$ws.Cells.Item(2, 12).Value = 1.0767494323679
$ws.Cells.Item(2, 13).Value = 1.088842776693001
$ws.Cells.Item(2, 14).Value = 1.079482627983157
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.074254730500714
$ws.Cells.Item(3, 4).Value = 1.07714950861477
$ws.Cells.Item(3, 5).Value = 1.075143636495122
$ws.Cells.Item(3, 6).Value = 1.08729335775987
$ws.Cells.Item(3, 9).Value = 1.058486924508238
$ws.Cells.Item(3, 10).Value = 1.078828530708353
$ws.Cells.Item(3, 11).Value = 1.079651204156306
$ws.Cells.Item(3, 12).Value = 1.077650248025135
$ws.Cells.Item(3, 13).Value = 1.089770491273244
$ws.Cells.Item(3, 14).Value = 1.080360590965321
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.075043614630836
$ws.Cells.Item(4, 4).Value = 1.077788782794776
$ws.Cells.Item(4, 5).Value = 1.075845869750061
$ws.Cells.Item(4, 6).Value = 1.088011310223659
$ws.Cells.Item(4, 9).Value = 1.058740790226224
$ws.Cells.Item(4, 10).Value = 1.079395072605341
$ws.Cells.Item(4, 11).Value = 1.080171467164715
$ws.Cells.Item(4, 12).Value = 1.078233075864264
$ws.Cells.Item(4, 13).Value = 1.090370488672984
$ws.Cells.Item(4, 14).Value = 1.080927937416786
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.075375219021572
$ws.Cells.Item(5, 4).Value = 1.078057497332646
$ws.Cells.Item(5, 5).Value = 1.076141218618232
$ws.Cells.Item(5, 6).Value = 1.088313210547838
$ws.Cells.Item(5, 9).Value = 1.058847212467398
$ws.Cells.Item(5, 10).Value = 1.079633067400878
$ws.Cells.Item(5, 11).Value = 1.080390001082363
$ws.Cells.Item(5, 12).Value = 1.078478082733782
$ws.Cells.Item(5, 13).Value = 1.090622657349216
$ws.Cells.Item(5, 14).Value = 1.08116627019223
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.075430894413834
$ws.Cells.Item(6, 4).Value = 1.078102613574284
$ws.Cells.Item(6, 5).Value = 1.076190816628812
$ws.Cells.Item(6, 6).Value = 1.088363905246319
$ws.Cells.Item(6, 9).Value = 1.058865063479225
$ws.Cells.Item(6, 10).Value = 1.079673017279554
$ws.Cells.Item(6, 11).Value = 1.080426683081174
$ws.Cells.Item(6, 12).Value = 1.078519219659813
$ws.Cells.Item(6, 13).Value = 1.090664993479383
$ws.Cells.Item(6, 14).Value = 1.081206276804317
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.075048045709847
$ws.Cells.Item(7, 4).Value = 1.077792373513522
$ws.Cells.Item(7, 5).Value = 1.07584981570472
$ws.Cells.Item(7, 6).Value = 1.088015343942501
$ws.Cells.Item(7, 9).Value = 1.058742213435124
$ws.Cells.Item(7, 10).Value = 1.079398253407297
$ws.Cells.Item(7, 11).Value = 1.080174387950289
$ws.Cells.Item(7, 12).Value = 1.078236349712323
$ws.Cells.Item(7, 13).Value = 1.090373858438352
$ws.Cells.Item(7, 14).Value = 1.080931122735846
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.073447435604095
$ws.Cells.Item(8, 4).Value = 1.076495310743158
$ws.Cells.Item(8, 5).Value = 1.074425589782329
$ws.Cells.Item(8, 6).Value = 1.0865590334032
$ws.Cells.Item(8, 9).Value = 1.058226148537247
$ws.Cells.Item(8, 10).Value = 1.078248260208713
$ws.Cells.Item(8, 11).Value = 1.079118266555583
$ws.Cells.Item(8, 12).Value = 1.077053879951826
$ws.Cells.Item(8, 13).Value = 1.08915636383729
$ws.Cells.Item(8, 14).Value = 1.079779496415014
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.070625469767532
$ws.Cells.Item(9, 4).Value = 1.074208508756972
$ws.Cells.Item(9, 5).Value = 1.071920048001175
$ws.Cells.Item(9, 6).Value = 1.083995117743192
$ws.Cells.Item(9, 9).Value = 1.057307054468826
$ws.Cells.Item(9, 10).Value = 1.076216015406083
$ws.Cells.Item(9, 11).Value = 1.077251284145195
$ws.Cells.Item(9, 12).Value = 1.074969736525737
$ws.Cells.Item(9, 13).Value = 1.087008687670265
$ws.Cells.Item(9, 14).Value = 1.077744365591662
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.068743054100053
$ws.Cells.Item(10, 4).Value = 1.072683116125223
$ws.Cells.Item(10, 5).Value = 1.070252431141156
$ws.Cells.Item(10, 6).Value = 1.08228733674406
$ws.Cells.Item(10, 9).Value = 1.056687761691079
$ws.Cells.Item(10, 10).Value = 1.074857219799007
$ws.Cells.Item(10, 11).Value = 1.07600257424424
$ws.Cells.Item(10, 12).Value = 1.073579953793238
$ws.Cells.Item(10, 13).Value = 1.085575318288789
$ws.Cells.Item(10, 14).Value = 1.076383640338968
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.067927660941552
$ws.Cells.Item(11, 4).Value = 1.072022390287662
$ws.Cells.Item(11, 5).Value = 1.06953097799046
$ws.Cells.Item(11, 6).Value = 1.081548193052232
$ws.Cells.Item(11, 9).Value = 1.056418037881187
$ws.Cells.Item(11, 10).Value = 1.074267891720237
$ws.Cells.Item(11, 11).Value = 1.075460896176905
$ws.Cells.Item(11, 12).Value = 1.072978071010695
$ws.Cells.Item(11, 13).Value = 1.084954267343709
$ws.Cells.Item(11, 14).Value = 1.075793475346729
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.067624741143218
$ws.Cells.Item(12, 4).Value = 1.071776933119697
$ws.Cells.Item(12, 5).Value = 1.06926309296134
$ws.Cells.Item(12, 6).Value = 1.081273691757606
$ws.Cells.Item(12, 9).Value = 1.056317614559313
$ws.Cells.Item(12, 10).Value = 1.074048843411197
$ws.Cells.Item(12, 11).Value = 1.07525954477267
$ws.Cells.Item(12, 12).Value = 1.072754489520604
$ws.Cells.Item(12, 13).Value = 1.084723521420071
$ws.Cells.Item(12, 14).Value = 1.075574115963966
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.067689720606482
$ws.Cells.Item(13, 4).Value = 1.071829586078038
$ws.Cells.Item(13, 5).Value = 1.069320550931911
$ws.Cells.Item(13, 6).Value = 1.081332570989633
$ws.Cells.Item(13, 9).Value = 1.056339166387245
$ws.Cells.Item(13, 10).Value = 1.074095836627589
$ws.Cells.Item(13, 11).Value = 1.075302742048951
$ws.Cells.Item(13, 12).Value = 1.072802449219293
$ws.Cells.Item(13, 13).Value = 1.084773019932531
$ws.Cells.Item(13, 14).Value = 1.075621175916115
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.067902622457253
$ws.Cells.Item(14, 4).Value = 1.072002101424535
$ws.Cells.Item(14, 5).Value = 1.069508832600973
$ws.Cells.Item(14, 6).Value = 1.08152550168396
$ws.Cells.Item(14, 9).Value = 1.056409741673748
$ws.Cells.Item(14, 10).Value = 1.074249788096817
$ws.Cells.Item(14, 11).Value = 1.07544425543332
$ws.Cells.Item(14, 12).Value = 1.072959590013997
$ws.Cells.Item(14, 13).Value = 1.084935195042291
$ws.Cells.Item(14, 14).Value = 1.075775346014088
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.068033792090373
$ws.Cells.Item(15, 4).Value = 1.072108389264838
$ws.Cells.Item(15, 5).Value = 1.069624851689956
$ws.Cells.Item(15, 6).Value = 1.081644379192221
$ws.Cells.Item(15, 9).Value = 1.056453194160774
$ws.Cells.Item(15, 10).Value = 1.074344623343001
$ws.Cells.Item(15, 11).Value = 1.075531426850166
$ws.Cells.Item(15, 12).Value = 1.073056407577412
$ws.Cells.Item(15, 13).Value = 1.085035108506681
$ws.Cells.Item(15, 14).Value = 1.075870315937199
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.068797162661499
$ws.Cells.Item(16, 4).Value = 1.072726961585995
$ws.Cells.Item(16, 5).Value = 1.070300324945467
$ws.Cells.Item(16, 6).Value = 1.082336398285353
$ws.Cells.Item(16, 9).Value = 1.056705629326058
$ws.Cells.Item(16, 10).Value = 1.074896311212776
$ws.Cells.Item(16, 11).Value = 1.076038502929201
$ws.Cells.Item(16, 12).Value = 1.073619896649158
$ws.Cells.Item(16, 13).Value = 1.085616527040987
$ws.Cells.Item(16, 14).Value = 1.076422787267029
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.06927592439236
$ws.Cells.Item(17, 4).Value = 1.073114915846272
$ws.Cells.Item(17, 5).Value = 1.070724201205514
$ws.Cells.Item(17, 6).Value = 1.082770573024155
$ws.Cells.Item(17, 9).Value = 1.056863555546274
$ws.Cells.Item(17, 10).Value = 1.07524211233695
$ws.Cells.Item(17, 11).Value = 1.07635631565766
$ws.Cells.Item(17, 12).Value = 1.073973331773677
$ws.Cells.Item(17, 13).Value = 1.085981129743032
$ws.Cells.Item(17, 14).Value = 1.076769079468463
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.069555149719893
$ws.Cells.Item(18, 4).Value = 1.073341181771594
$ws.Cells.Item(18, 5).Value = 1.07097150246908
$ws.Cells.Item(18, 6).Value = 1.083023852364814
$ws.Cells.Item(18, 9).Value = 1.056955520189782
$ws.Cells.Item(18, 10).Value = 1.075443719685732
$ws.Cells.Item(18, 11).Value = 1.076541596112096
$ws.Cells.Item(18, 12).Value = 1.074179475128469
$ws.Cells.Item(18, 13).Value = 1.086193758522811
$ws.Cells.Item(18, 14).Value = 1.076970973122804
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.069650353608021
$ws.Cells.Item(19, 4).Value = 1.073418329033669
$ws.Cells.Item(19, 5).Value = 1.071055836247937
$ws.Cells.Item(19, 6).Value = 1.083110219627609
$ws.Cells.Item(19, 9).Value = 1.056986852154751
$ws.Cells.Item(19, 10).Value = 1.075512446914568
$ws.Cells.Item(19, 11).Value = 1.076604755989862
$ws.Cells.Item(19, 12).Value = 1.074249763117498
$ws.Cells.Item(19, 13).Value = 1.086266253058074
$ws.Cells.Item(19, 14).Value = 1.077039797952189
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.06922456073279
$ws.Cells.Item(20, 4).Value = 1.073073294229903
$ws.Cells.Item(20, 5).Value = 1.070678716948481
$ws.Cells.Item(20, 6).Value = 1.082723986823975
$ws.Cells.Item(20, 9).Value = 1.056846627186123
$ws.Cells.Item(20, 10).Value = 1.075205020736921
$ws.Cells.Item(20, 11).Value = 1.076322227125379
$ws.Cells.Item(20, 12).Value = 1.073935412505476
$ws.Cells.Item(20, 13).Value = 1.085942015242211
$ws.Cells.Item(20, 14).Value = 1.076731935194107
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.067839929469147
$ws.Cells.Item(21, 4).Value = 1.07195130090626
$ws.Cells.Item(21, 5).Value = 1.069453385731502
$ws.Cells.Item(21, 6).Value = 1.081468687042544
$ws.Cells.Item(21, 9).Value = 1.05638896552762
$ws.Cells.Item(21, 10).Value = 1.074204457256556
$ws.Cells.Item(21, 11).Value = 1.075402587364163
$ws.Cells.Item(21, 12).Value = 1.07291331639844
$ws.Cells.Item(21, 13).Value = 1.084887440177337
$ws.Cells.Item(21, 14).Value = 1.075729950798834
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.06696908686871
$ws.Cells.Item(22, 4).Value = 1.071245660957096
$ws.Cells.Item(22, 5).Value = 1.068683519436334
$ws.Cells.Item(22, 6).Value = 1.080679716617019
$ws.Cells.Item(22, 9).Value = 1.056099850386753
$ws.Cells.Item(22, 10).Value = 1.073574520123231
$ws.Cells.Item(22, 11).Value = 1.074823516241888
$ws.Cells.Item(22, 12).Value = 1.072270594388104
$ws.Cells.Item(22, 13).Value = 1.084224039640109
$ws.Cells.Item(22, 14).Value = 1.075099119082524
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.067430763330135
$ws.Cells.Item(23, 4).Value = 1.071619753170526
$ws.Cells.Item(23, 5).Value = 1.069091588426084
$ws.Cells.Item(23, 6).Value = 1.08109793796082
$ws.Cells.Item(23, 9).Value = 1.056253245396104
$ws.Cells.Item(23, 10).Value = 1.073908542072699
$ws.Cells.Item(23, 11).Value = 1.0751305743956
$ws.Cells.Item(23, 12).Value = 1.072611322232182
$ws.Cells.Item(23, 13).Value = 1.08457575406658
$ws.Cells.Item(23, 14).Value = 1.075433615381474
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.069247769847101
$ws.Cells.Item(24, 4).Value = 1.073092101312503
$ws.Cells.Item(24, 5).Value = 1.070699269136819
$ws.Cells.Item(24, 6).Value = 1.082745037022114
$ws.Cells.Item(24, 9).Value = 1.056854276850175
$ws.Cells.Item(24, 10).Value = 1.075221781120583
$ws.Cells.Item(24, 11).Value = 1.076337630557369
$ws.Cells.Item(24, 12).Value = 1.073952546619672
$ws.Cells.Item(24, 13).Value = 1.085959689518421
$ws.Cells.Item(24, 14).Value = 1.076748719379436
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.071355202606751
$ws.Cells.Item(25, 4).Value = 1.074799849924654
$ws.Cells.Item(25, 5).Value = 1.072567304701992
$ws.Cells.Item(25, 6).Value = 1.084657684520156
$ws.Cells.Item(25, 9).Value = 1.057545817180882
$ws.Cells.Item(25, 10).Value = 1.07674209402252
$ws.Cells.Item(25, 11).Value = 1.077734655134437
$ws.Cells.Item(25, 12).Value = 1.075508597664264
$ws.Cells.Item(25, 13).Value = 1.08756418975865
$ws.Cells.Item(25, 14).Value = 1.078271191300076
